$d = $word.ActiveDocument

# Update the header date line
$d.Content.Find.Execute("2023-11-20 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-11-21 Tuesday", 2)

# Update the division problems in the table (row 1,5,9,13,17 hold the
# problems; the rows in between are blank "work space" rows).
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "98÷9="
$t.Cell(1, 2).Range.Text = "71÷8="
$t.Cell(1, 3).Range.Text = "15÷8="
$t.Cell(1, 4).Range.Text = "79÷7="
$t.Cell(1, 5).Range.Text = "28÷7="

$t.Cell(5, 1).Range.Text = "78÷8="
$t.Cell(5, 2).Range.Text = "30÷2="
$t.Cell(5, 3).Range.Text = "72÷7="
$t.Cell(5, 4).Range.Text = "87÷7="
$t.Cell(5, 5).Range.Text = "18÷4="

$t.Cell(9, 1).Range.Text = "33÷4="
$t.Cell(9, 2).Range.Text = "16÷4="
$t.Cell(9, 3).Range.Text = "94÷7="
$t.Cell(9, 4).Range.Text = "79÷4="
$t.Cell(9, 5).Range.Text = "86÷6="

$t.Cell(13, 1).Range.Text = "78÷2="
$t.Cell(13, 2).Range.Text = "99÷6="
$t.Cell(13, 3).Range.Text = "62÷8="
$t.Cell(13, 4).Range.Text = "17÷9="
$t.Cell(13, 5).Range.Text = "74÷8="

$t.Cell(17, 1).Range.Text = "53÷2="
$t.Cell(17, 2).Range.Text = "20÷6="
$t.Cell(17, 3).Range.Text = "99÷3="
$t.Cell(17, 4).Range.Text = "86÷4="
$t.Cell(17, 5).Range.Text = "20÷5="
